{"js": "// Replace the date line and every \"NN\u00d7NN=NNNN\" answer cell with the\n// values from the newer worksheet generation. Each old value below is\n// unique within the document, so a plain case-sensitive search finds\n// exactly one run to update.\nconst replacements = [\n  [\"2026-01-28 Wednesday\", \"2026-01-29 Thursday\"],\n  [\"88\u00d763=5544\", \"47\u00d756=2632\"],\n  [\"31\u00d782=2542\", \"65\u00d789=5785\"],\n  [\"92\u00d773=6716\", \"41\u00d713=533\"],\n  [\"33\u00d761=2013\", \"20\u00d748=960\"],\n  [\"65\u00d778=5070\", \"57\u00d741=2337\"],\n  [\"46\u00d711=506\", \"88\u00d785=7480\"],\n  [\"17\u00d719=323\", \"25\u00d713=325\"],\n  [\"53\u00d799=5247\", \"99\u00d753=5247\"],\n  [\"36\u00d712=432\", \"59\u00d779=4661\"],\n  [\"42\u00d780=3360\", \"95\u00d768=6460\"],\n  [\"34\u00d748=1632\", \"64\u00d713=832\"],\n  [\"32\u00d784=2688\", \"44\u00d793=4092\"],\n  [\"59\u00d773=4307\", \"16\u00d749=784\"],\n  [\"30\u00d761=1830\", \"45\u00d777=3465\"],\n  [\"20\u00d782=1640\", \"56\u00d730=1680\"],\n  [\"88\u00d774=6512\", \"68\u00d796=6528\"],\n  [\"90\u00d754=4860\", \"63\u00d763=3969\"],\n  [\"50\u00d721=1050\", \"14\u00d753=742\"],\n  [\"54\u00d752=2808\", \"71\u00d759=4189\"],\n  [\"59\u00d725=1475\", \"41\u00d718=738\"],\n  [\"79\u00d776=6004\", \"29\u00d762=1798\"],\n  [\"97\u00d714=1358\", \"87\u00d783=7221\"],\n  [\"88\u00d786=7568\", \"90\u00d757=5130\"],\n  [\"96\u00d757=5472\", \"27\u00d759=1593\"],\n  [\"91\u00d746=4186\", \"52\u00d714=728\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and every \"NN\u00d7NN=NNNN\" answer cell with the\n# values from the newer worksheet generation. Each old value is unique\n# within the document, so Find/Replace (ReplaceAll) touches exactly one\n# run per pair without disturbing its surrounding formatting.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-28 Wednesday\", \"2026-01-29 Thursday\"),\n    @(\"88\u00d763=5544\", \"47\u00d756=2632\"),\n    @(\"31\u00d782=2542\", \"65\u00d789=5785\"),\n    @(\"92\u00d773=6716\", \"41\u00d713=533\"),\n    @(\"33\u00d761=2013\", \"20\u00d748=960\"),\n    @(\"65\u00d778=5070\", \"57\u00d741=2337\"),\n    @(\"46\u00d711=506\", \"88\u00d785=7480\"),\n    @(\"17\u00d719=323\", \"25\u00d713=325\"),\n    @(\"53\u00d799=5247\", \"99\u00d753=5247\"),\n    @(\"36\u00d712=432\", \"59\u00d779=4661\"),\n    @(\"42\u00d780=3360\", \"95\u00d768=6460\"),\n    @(\"34\u00d748=1632\", \"64\u00d713=832\"),\n    @(\"32\u00d784=2688\", \"44\u00d793=4092\"),\n    @(\"59\u00d773=4307\", \"16\u00d749=784\"),\n    @(\"30\u00d761=1830\", \"45\u00d777=3465\"),\n    @(\"20\u00d782=1640\", \"56\u00d730=1680\"),\n    @(\"88\u00d774=6512\", \"68\u00d796=6528\"),\n    @(\"90\u00d754=4860\", \"63\u00d763=3969\"),\n    @(\"50\u00d721=1050\", \"14\u00d753=742\"),\n    @(\"54\u00d752=2808\", \"71\u00d759=4189\"),\n    @(\"59\u00d725=1475\", \"41\u00d718=738\"),\n    @(\"79\u00d776=6004\", \"29\u00d762=1798\"),\n    @(\"97\u00d714=1358\", \"87\u00d783=7221\"),\n    @(\"88\u00d786=7568\", \"90\u00d757=5130\"),\n    @(\"96\u00d757=5472\", \"27\u00d759=1593\"),\n    @(\"91\u00d746=4186\", \"52\u00d714=728\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2)\n}\n"}
